$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(5, 6).Value = 8993
$ws.Cells.Item(6, 6).Value = 250
$ws.Cells.Item(7, 6).Value = 109
$ws.Cells.Item(8, 6).Value = 1259
$ws.Cells.Item(10, 6).Value = 572
$ws.Cells.Item(11, 6).Value = 597
$ws.Cells.Item(13, 6).Value = 130
$ws.Cells.Item(14, 6).Value = 286
$ws.Cells.Item(16, 6).Value = 51
$ws.Cells.Item(17, 6).Value = 1467
$ws.Cells.Item(18, 6).Value = 1315
$ws.Cells.Item(21, 6).Value = 1355
$ws.Cells.Item(22, 6).Value = 74
$ws.Cells.Item(23, 6).Value = 223
$ws.Cells.Item(25, 6).Value = 81
$ws.Cells.Item(26, 6).Value = 46
$ws.Cells.Item(28, 6).Value = 291
$ws.Cells.Item(29, 6).Value = 291
$ws.Cells.Item(30, 6).Value = 1062
$ws.Cells.Item(31, 6).Value = 9
$ws.Cells.Item(33, 6).Value = 221
$ws.Cells.Item(34, 6).Value = 189
$ws.Cells.Item(37, 6).Value = 605
$ws.Cells.Item(38, 6).Value = 428
$ws.Cells.Item(42, 6).Value = 4
$ws.Cells.Item(44, 6).Value = 1223
$ws.Cells.Item(46, 6).Value = 202
$ws.Cells.Item(47, 6).Value = 42
$ws.Cells.Item(48, 6).Value = 42

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 13
$ws.Cells.Item(19, 6).Value = 8
$ws.Cells.Item(20, 6).Value = 69
$ws.Cells.Item(23, 6).Value = 927
$ws.Cells.Item(26, 6).Value = 216
$ws.Cells.Item(29, 6).Value = 196
$ws.Cells.Item(31, 6).Value = 147

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(5, 6).Value = 743
$ws.Cells.Item(6, 6).Value = 285
$ws.Cells.Item(7, 6).Value = 141
$ws.Cells.Item(8, 6).Value = 2023
$ws.Cells.Item(9, 6).Value = 3045

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(5, 6).Value = 743
$ws.Cells.Item(6, 6).Value = 8993
$ws.Cells.Item(7, 6).Value = 285
$ws.Cells.Item(8, 6).Value = 141
$ws.Cells.Item(9, 6).Value = 13
$ws.Cells.Item(10, 6).Value = 250
$ws.Cells.Item(11, 6).Value = 2023
$ws.Cells.Item(12, 6).Value = 3045
$ws.Cells.Item(14, 6).Value = 1259
$ws.Cells.Item(16, 6).Value = 572
$ws.Cells.Item(17, 6).Value = 597
$ws.Cells.Item(18, 6).Value = 286
$ws.Cells.Item(19, 6).Value = 51
$ws.Cells.Item(20, 6).Value = 1467
$ws.Cells.Item(21, 6).Value = 1315
$ws.Cells.Item(23, 6).Value = 1355
$ws.Cells.Item(24, 6).Value = 74
$ws.Cells.Item(26, 6).Value = 81
$ws.Cells.Item(27, 6).Value = 46
$ws.Cells.Item(28, 6).Value = 291
$ws.Cells.Item(29, 6).Value = 1062
$ws.Cells.Item(31, 6).Value = 69
$ws.Cells.Item(32, 6).Value = 221
$ws.Cells.Item(34, 6).Value = 216
$ws.Cells.Item(37, 6).Value = 605
$ws.Cells.Item(40, 6).Value = 196
$ws.Cells.Item(41, 6).Value = 147
$ws.Cells.Item(46, 6).Value = 202
$ws.Cells.Item(49, 6).Value = 42
